$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.346.13"
$ws.Range("E2").Value = "  +2.25%  "

# Row 3
$ws.Range("D3").Value = "3.915.81"
$ws.Range("E3").Value = "  +3.92%  "

# Row 4
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").Value = "'471.27"
$ws.Range("E5").Value = "  +10.41%  "

# Row 6
$ws.Range("D6").Value = "'144.60"
$ws.Range("E6").Value = "  +4.70%  "

# Row 7
$ws.Range("D7").Value = "'0.624"
$ws.Range("E7").Value = "  +1.02%  "

# Row 8
$ws.Range("E8").Value = "  -0.09%  "

# Row 9
$ws.Range("D9").Value = "'0.736"
$ws.Range("E9").Value = "  +1.73%  "

# Row 10
$ws.Range("D10").Value = "'0.164"
$ws.Range("E10").Value = "  +10.50%  "

# Row 11
$ws.Range("D11").Value = "'0.0000339"
$ws.Range("E11").Value = "  +11.61%  "

# Row 12
$ws.Range("D12").Value = "'43.35"
$ws.Range("E12").Value = "  +2.54%  "

# Row 13
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'10.40"
$ws.Range("E13").Value = "  +0.75%  "

# Row 14
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "4.527.06"
$ws.Range("E14").Value = "  +3.56%  "

# Row 15
$ws.Range("D15").Value = "'15.11"
$ws.Range("E15").Value = "  +1.10%  "

# Row 16
$ws.Range("D16").Value = "3.889.95"
$ws.Range("E16").Value = "  +2.68%  "

# Row 17
$ws.Range("E17").Value = "  -0.14%  "

# Row 18
$ws.Range("D18").Value = "'19.89"
$ws.Range("E18").Value = "  +0.73%  "

# Row 19
$ws.Range("E19").Value = "  +4.42%  "

# Row 20
$ws.Range("D20").Value = "67.618.35"
$ws.Range("E20").Value = "  +2.56%  "

# Row 21
$ws.Range("D21").Value = "'433.00"
$ws.Range("E21").Value = "  +7.57%  "

# Row 22
$ws.Range("B22").Value = "InternetComputer(DFINITY)"
$ws.Range("C22").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D22").Value = "'14.63"
$ws.Range("E22").Value = "  -0.99%  "

# Row 23
$ws.Range("B23").Value = "ImmutableX"
$ws.Range("C23").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D23").Value = "'3.36"
$ws.Range("E23").Value = "  +3.06%  "

# Row 24
$ws.Range("D24").Value = "'88.26"
$ws.Range("E24").Value = "  +4.97%  "

# Row 25
$ws.Range("D25").Value = "'3.57"
$ws.Range("E25").Value = "  +9.63%  "

# Row 26
$ws.Range("D26").Value = "'38.62"
$ws.Range("E26").Value = "  +6.15%  "

# Row 27
$ws.Range("D27").Value = "'5.75"
$ws.Range("E27").Value = "  +5.02%  "

# Row 28
$ws.Range("D28").Value = "'10.09"
$ws.Range("E28").Value = "  +3.45%  "

# Row 29
$ws.Range("D29").Value = "'9.54"
$ws.Range("E29").Value = "  -4.87%  "

# Row 30
$ws.Range("D30").Value = "'726.07"
$ws.Range("E30").Value = "  +3.78%  "

# Row 31
$ws.Range("D31").Value = "'13.73"
$ws.Range("E31").Value = "  +0.83%  "

# Row 32
$ws.Range("E32").Value = "  -0.20%  "

# Row 33
$ws.Range("D33").Value = "'2.82"
$ws.Range("E33").Value = "  +2.01%  "

# Row 34
$ws.Range("D34").Value = "'43.20"
$ws.Range("E34").Value = "  +6.22%  "

# Row 35
$ws.Range("E35").Value = "  +5.33%  "

# Row 36
$ws.Range("D36").Value = "'57.53"
$ws.Range("E36").Value = "  +2.54%  "

# Row 37
$ws.Range("D37").Value = "0.0₃0801"
$ws.Range("E37").Value = "  +21.72%  "

# Row 38
$ws.Range("E38").Value = "  +0.19%  "

# Row 39
$ws.Range("D39").Value = "'5.40"
$ws.Range("E39").Value = "  -5.08%  "

# Row 40
$ws.Range("D40").Value = "'0.0476"
$ws.Range("E40").Value = "  +2.22%  "

# Row 41
$ws.Range("E41").Value = "  +4.19%  "

# Row 42
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").Value = "'0.141"
$ws.Range("E42").Value = "  +0.84%  "

# Row 43
$ws.Range("B43").Value = "TheGraph"
$ws.Range("C43").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D43").Value = "'0.336"
$ws.Range("E43").Value = "  +4.97%  "

# Row 44
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "  +0.05%  "

# Row 45
$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D45").Value = "'2.55"
$ws.Range("E45").Value = "  -6.56%  "

# Row 46
$ws.Range("E46").Value = "  +4.49%  "

# Row 47
$ws.Range("D47").Value = "'2.18"
$ws.Range("E47").Value = "  +6.64%  "

# Row 48
$ws.Range("D48").Value = "'3.40"
$ws.Range("E48").Value = "  +2.09%  "

# Row 49
$ws.Range("D49").Value = "'3.17"
$ws.Range("E49").Value = "  -0.79%  "

# Row 50
$ws.Range("D50").Value = "'145.34"
$ws.Range("E50").Value = "  +5.42%  "

# Row 51
$ws.Range("E51").Value = "  +4.59%  "
